$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to the Price column (D) while forcing text storage,
# matching the workbook's original inlineStr string cells (e.g. "534.82").
# A leading apostrophe tells Excel to treat the numeric-looking text as a
# literal string (quote-prefix) instead of silently converting it to a number.
function Set-PriceText($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "59.444.08"
$ws.Range("E2").Value = "  +2.79%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.172.80"
$ws.Range("E3").Value = "  +1.58%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-PriceText "D5" "534.68"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6 - Solana
Set-PriceText "D6" "144.19"
$ws.Range("E6").Value = "  +3.66%  "

# Row 7 - USDC
Set-PriceText "D7" "1.00"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - XRP
Set-PriceText "D8" "0.520"
$ws.Range("E8").Value = "  +5.80%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  -1.21%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.51%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +3.89%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PriceText "D12" "3.719.99"
$ws.Range("E12").Value = "  +1.54%  "

# Row 13 - TRON
Set-PriceText "D13" "0.139"
$ws.Range("E13").Value = "  -0.05%  "

# Row 14 - Avalanche
Set-PriceText "D14" "25.97"
$ws.Range("E14").Value = "  +0.37%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +3.37%  "

# Row 16 - WrappedBTC
Set-PriceText "D16" "59.467.69"
$ws.Range("E16").Value = "  +2.63%  "

# Row 17 - WrappedEther
Set-PriceText "D17" "3.210.19"
$ws.Range("E17").Value = "  +2.83%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +1.29%  "

# Row 19 - Chainlink
Set-PriceText "D19" "12.95"
$ws.Range("E19").Value = "  +0.94%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +0.45%  "

# Row 21 - BitcoinCash
Set-PriceText "D21" "377.78"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - Polygon
Set-PriceText "D23" "0.528"
$ws.Range("E23").Value = "  +3.88%  "

# Row 24 - Litecoin
Set-PriceText "D24" "70.02"

# Row 25 - was InternetComputer(DFINITY), now Kaspa
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-PriceText "D25" "0.170"
$ws.Range("E25").Value = "  +1.59%  "

# Row 26 - was Kaspa, now InternetComputer(DFINITY)
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-PriceText "D26" "8.80"
$ws.Range("E26").Value = "  +16.68%  "

# Row 27 - Binance-PegBSC-USD
Set-PriceText "D27" "0.998"

# Row 28 - PEPE
$ws.Range("E28").Value = "  +2.57%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +1.65%  "

# Row 30 - EthereumClassic
Set-PriceText "D30" "22.30"
$ws.Range("E30").Value = "  +3.40%  "

# Row 31 - RenderToken
Set-PriceText "D31" "6.13"
$ws.Range("E31").Value = "  -0.44%  "

# Row 32 - NEARProtocol
Set-PriceText "D32" "5.29"
$ws.Range("E32").Value = "  +2.05%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -0.04%  "

# Row 34 - Aptos
Set-PriceText "D34" "6.46"
$ws.Range("E34").Value = "  +4.54%  "

# Row 35 - Monero
Set-PriceText "D35" "156.42"
$ws.Range("E35").Value = "  -2.37%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +4.35%  "

# Row 37 - was Hedera, now Maker
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-PriceText "D37" "2.739.64"
$ws.Range("E37").Value = "  +6.60%  "

# Row 38 - was Maker, now Hedera
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-PriceText "D38" "0.0710"
$ws.Range("E38").Value = "  +5.68%  "

# Row 39 - EnergySwap
Set-PriceText "D39" "25.32"
$ws.Range("E39").Value = "  -1.24%  "

# Row 40 - Stacks
Set-PriceText "D40" "1.67"
$ws.Range("E40").Value = "  +1.75%  "

# Row 41 - Filecoin
Set-PriceText "D41" "4.28"
$ws.Range("E41").Value = "  +3.11%  "

# Row 42 - was OKB, now Mantle
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceText "D42" "0.723"
$ws.Range("E42").Value = "  +3.53%  "

# Row 43 - was Mantle, now OKB
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-PriceText "D43" "39.34"
$ws.Range("E43").Value = "  +3.13%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +7.25%  "

# Row 45 - was RenzoRestakedETH, now ONDO
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-PriceText "D45" "1.00"
$ws.Range("E45").Value = "  +2.24%  "

# Row 46 - was ONDO, now RenzoRestakedETH
$ws.Range("B46").Value = "RenzoRestakedETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-PriceText "D46" "3.215.32"
$ws.Range("E46").Value = "  +1.55%  "

# Row 47 - Cosmos
$ws.Range("E47").Value = "  +0.35%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +7.26%  "

# Row 49 - InjectiveProtocol
Set-PriceText "D49" "20.45"
$ws.Range("E49").Value = "  +2.97%  "

# Row 50 - FirstDigitalUSD
$ws.Range("E50").Value = "  -0.04%  "

# Row 51 - SuiNetwork
Set-PriceText "D51" "0.768"
$ws.Range("E51").Value = "  +2.26%  "
